$wb = $excel.ActiveWorkbook

# --- Generic: NrBuckets 8 -> 16 ---
$ws1 = $wb.Worksheets.Item("Generic")
$ws1.Cells.Item(4,2).Value = 16

# --- Productdata ---
$ws3 = $wb.Worksheets.Item("Productdata")
$ws3.Cells.Item(2,3).Value = 0
$ws3.Cells.Item(2,5).Value = 195.3772307692308
$ws3.Cells.Item(3,5).Value = 33.0097373076923
$ws3.Cells.Item(4,5).Value = 31.93882692307693
$ws3.Cells.Item(5,5).Value = 75.440715
$ws3.Cells.Item(6,5).Value = 85.2724380769231
$ws3.Cells.Item(7,3).Value = 3425
$ws3.Cells.Item(7,5).Value = 33.0097373076923
$ws3.Cells.Item(8,3).Value = 2868
$ws3.Cells.Item(8,5).Value = 31.93882692307693
$ws3.Cells.Item(9,3).Value = 14071
$ws3.Cells.Item(9,5).Value = 150.88143
$ws3.Cells.Item(10,3).Value = 1522
$ws3.Cells.Item(10,5).Value = 9.831723076923078
$ws3.Cells.Item(11,3).Value = 8839
$ws3.Cells.Item(11,5).Value = 28.72970653846154
$ws3.Cells.Item(12,3).Value = 6429
$ws3.Cells.Item(12,5).Value = 27.79541153846154
$ws3.Cells.Item(13,3).Value = 12273
$ws3.Cells.Item(13,5).Value = 65.637915
$ws3.Cells.Item(14,3).Value = 14682
$ws3.Cells.Item(14,5).Value = 74.19209961538462

# --- ForecastedAverageDemand: update existing rows, then append rows 10-17 ---
$ws4 = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws4.Cells.Item(2,8).Value = 1872
$ws4.Cells.Item(2,9).Value = 12113
$ws4.Cells.Item(2,10).Value = 0
$ws4.Cells.Item(3,7).Value = 3385
$ws4.Cells.Item(3,8).Value = 973
$ws4.Cells.Item(3,9).Value = 1859
$ws4.Cells.Item(3,10).Value = 1504
$ws4.Cells.Item(4,7).Value = 8690
$ws4.Cells.Item(4,8).Value = 1095
$ws4.Cells.Item(4,9).Value = 2531
$ws4.Cells.Item(4,10).Value = 2369
$ws4.Cells.Item(5,7).Value = 0
$ws4.Cells.Item(5,8).Value = 5202
$ws4.Cells.Item(5,9).Value = 9492
$ws4.Cells.Item(5,10).Value = 0
$ws4.Cells.Item(6,7).Value = 1051
$ws4.Cells.Item(6,8).Value = 7632
$ws4.Cells.Item(6,9).Value = 0
$ws4.Cells.Item(6,10).Value = 3370
$ws4.Cells.Item(7,7).Value = 4565
$ws4.Cells.Item(7,8).Value = 6607
$ws4.Cells.Item(7,9).Value = 6646
$ws4.Cells.Item(7,10).Value = 0
$ws4.Cells.Item(8,7).Value = 0
$ws4.Cells.Item(8,8).Value = 0
$ws4.Cells.Item(8,9).Value = 0
$ws4.Cells.Item(8,10).Value = 0
$ws4.Cells.Item(9,7).Value = 4229
$ws4.Cells.Item(9,8).Value = 1784
$ws4.Cells.Item(9,9).Value = 4348
$ws4.Range("A9:N9").Copy()
$ws4.Range("A10:N17").PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false
$ws4.Cells.Item(10,1).Value = 8
$ws4.Cells.Item(10,2).Value = 0
$ws4.Cells.Item(10,3).Value = 0
$ws4.Cells.Item(10,4).Value = 0
$ws4.Cells.Item(10,5).Value = 0
$ws4.Cells.Item(10,6).Value = 0
$ws4.Cells.Item(10,7).Value = 1681
$ws4.Cells.Item(10,8).Value = 0
$ws4.Cells.Item(10,9).Value = 9137
$ws4.Cells.Item(10,10).Value = 0
$ws4.Cells.Item(10,11).Value = 0
$ws4.Cells.Item(10,12).Value = 0
$ws4.Cells.Item(10,13).Value = 0
$ws4.Cells.Item(10,14).Value = 0
$ws4.Cells.Item(11,1).Value = 9
$ws4.Cells.Item(11,2).Value = 0
$ws4.Cells.Item(11,3).Value = 0
$ws4.Cells.Item(11,4).Value = 0
$ws4.Cells.Item(11,5).Value = 0
$ws4.Cells.Item(11,6).Value = 0
$ws4.Cells.Item(11,7).Value = 2472
$ws4.Cells.Item(11,8).Value = 2402
$ws4.Cells.Item(11,9).Value = 16920
$ws4.Cells.Item(11,10).Value = 766
$ws4.Cells.Item(11,11).Value = 0
$ws4.Cells.Item(11,12).Value = 0
$ws4.Cells.Item(11,13).Value = 0
$ws4.Cells.Item(11,14).Value = 0
$ws4.Cells.Item(12,1).Value = 10
$ws4.Cells.Item(12,2).Value = 0
$ws4.Cells.Item(12,3).Value = 0
$ws4.Cells.Item(12,4).Value = 0
$ws4.Cells.Item(12,5).Value = 0
$ws4.Cells.Item(12,6).Value = 0
$ws4.Cells.Item(12,7).Value = 8670
$ws4.Cells.Item(12,8).Value = 7577
$ws4.Cells.Item(12,9).Value = 1230
$ws4.Cells.Item(12,10).Value = 1596
$ws4.Cells.Item(12,11).Value = 0
$ws4.Cells.Item(12,12).Value = 0
$ws4.Cells.Item(12,13).Value = 0
$ws4.Cells.Item(12,14).Value = 0
$ws4.Cells.Item(13,1).Value = 11
$ws4.Cells.Item(13,2).Value = 0
$ws4.Cells.Item(13,3).Value = 0
$ws4.Cells.Item(13,4).Value = 0
$ws4.Cells.Item(13,5).Value = 0
$ws4.Cells.Item(13,6).Value = 0
$ws4.Cells.Item(13,7).Value = 7635
$ws4.Cells.Item(13,8).Value = 6153
$ws4.Cells.Item(13,9).Value = 13606
$ws4.Cells.Item(13,10).Value = 0
$ws4.Cells.Item(13,11).Value = 0
$ws4.Cells.Item(13,12).Value = 0
$ws4.Cells.Item(13,13).Value = 0
$ws4.Cells.Item(13,14).Value = 0
$ws4.Cells.Item(14,1).Value = 12
$ws4.Cells.Item(14,2).Value = 0
$ws4.Cells.Item(14,3).Value = 0
$ws4.Cells.Item(14,4).Value = 0
$ws4.Cells.Item(14,5).Value = 0
$ws4.Cells.Item(14,6).Value = 0
$ws4.Cells.Item(14,7).Value = 712
$ws4.Cells.Item(14,8).Value = 0
$ws4.Cells.Item(14,9).Value = 11239
$ws4.Cells.Item(14,10).Value = 1480
$ws4.Cells.Item(14,11).Value = 0
$ws4.Cells.Item(14,12).Value = 0
$ws4.Cells.Item(14,13).Value = 0
$ws4.Cells.Item(14,14).Value = 0
$ws4.Cells.Item(15,1).Value = 13
$ws4.Cells.Item(15,2).Value = 0
$ws4.Cells.Item(15,3).Value = 0
$ws4.Cells.Item(15,4).Value = 0
$ws4.Cells.Item(15,5).Value = 0
$ws4.Cells.Item(15,6).Value = 0
$ws4.Cells.Item(15,7).Value = 3277
$ws4.Cells.Item(15,8).Value = 0
$ws4.Cells.Item(15,9).Value = 2804
$ws4.Cells.Item(15,10).Value = 2755
$ws4.Cells.Item(15,11).Value = 0
$ws4.Cells.Item(15,12).Value = 0
$ws4.Cells.Item(15,13).Value = 0
$ws4.Cells.Item(15,14).Value = 0
$ws4.Cells.Item(16,1).Value = 14
$ws4.Cells.Item(16,2).Value = 0
$ws4.Cells.Item(16,3).Value = 0
$ws4.Cells.Item(16,4).Value = 0
$ws4.Cells.Item(16,5).Value = 0
$ws4.Cells.Item(16,6).Value = 0
$ws4.Cells.Item(16,7).Value = 0
$ws4.Cells.Item(16,8).Value = 48
$ws4.Cells.Item(16,9).Value = 14272
$ws4.Cells.Item(16,10).Value = 0
$ws4.Cells.Item(16,11).Value = 0
$ws4.Cells.Item(16,12).Value = 0
$ws4.Cells.Item(16,13).Value = 0
$ws4.Cells.Item(16,14).Value = 0
$ws4.Cells.Item(17,1).Value = 15
$ws4.Cells.Item(17,2).Value = 0
$ws4.Cells.Item(17,3).Value = 0
$ws4.Cells.Item(17,4).Value = 0
$ws4.Cells.Item(17,5).Value = 0
$ws4.Cells.Item(17,6).Value = 0
$ws4.Cells.Item(17,7).Value = 0
$ws4.Cells.Item(17,8).Value = 3542
$ws4.Cells.Item(17,9).Value = 0
$ws4.Cells.Item(17,10).Value = 0
$ws4.Cells.Item(17,11).Value = 0
$ws4.Cells.Item(17,12).Value = 0
$ws4.Cells.Item(17,13).Value = 0
$ws4.Cells.Item(17,14).Value = 0

# --- ForcastedStandardDeviation: update existing rows, then append rows 10-17 ---
$ws5 = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws5.Cells.Item(2,8).Value = 46.79999999999999
$ws5.Cells.Item(2,9).Value = 302.8249999999999
$ws5.Cells.Item(2,10).Value = 0
$ws5.Cells.Item(3,7).Value = 160.7875
$ws5.Cells.Item(3,8).Value = 46.21749999999999
$ws5.Cells.Item(3,9).Value = 88.30249999999998
$ws5.Cells.Item(3,10).Value = 71.43999999999998
$ws5.Cells.Item(4,7).Value = 588.7474999999998
$ws5.Cells.Item(4,8).Value = 74.18624999999997
$ws5.Cells.Item(4,9).Value = 171.4752499999999
$ws5.Cells.Item(4,10).Value = 160.4997499999999
$ws5.Cells.Item(5,7).Value = 0
$ws5.Cells.Item(5,8).Value = 447.24195
$ws5.Cells.Item(5,9).Value = 816.0747
$ws5.Cells.Item(5,10).Value = 0
$ws5.Cells.Item(6,7).Value = 107.5987525
$ws5.Cells.Item(6,8).Value = 781.3450799999998
$ws5.Cells.Item(6,9).Value = 0
$ws5.Cells.Item(6,10).Value = 345.012175
$ws5.Cells.Item(7,7).Value = 534.74295875
$ws5.Cells.Item(7,8).Value = 773.94232825
$ws5.Cells.Item(7,9).Value = 778.5107784999999
$ws5.Cells.Item(7,10).Value = 0
$ws5.Cells.Item(8,7).Value = 0
$ws5.Cells.Item(8,8).Value = 0
$ws5.Cells.Item(8,9).Value = 0
$ws5.Cells.Item(8,10).Value = 0
$ws5.Cells.Item(9,7).Value = 602.1385422274999
$ws5.Cells.Item(9,8).Value = 254.01162434
$ws5.Cells.Item(9,9).Value = 619.0821427299999
$ws5.Range("A9:N9").Copy()
$ws5.Range("A10:N17").PasteSpecial(-4122)
$ws5.Application.CutCopyMode = $false
$ws5.Cells.Item(10,1).Value = 8
$ws5.Cells.Item(10,2).Value = 0
$ws5.Cells.Item(10,3).Value = 0
$ws5.Cells.Item(10,4).Value = 0
$ws5.Cells.Item(10,5).Value = 0
$ws5.Cells.Item(10,6).Value = 0
$ws5.Cells.Item(10,7).Value = 257.4365394977499
$ws5.Cells.Item(10,8).Value = 0
$ws5.Cells.Item(10,9).Value = 1399.28474800175
$ws5.Cells.Item(10,10).Value = 0
$ws5.Cells.Item(10,11).Value = 0
$ws5.Cells.Item(10,12).Value = 0
$ws5.Cells.Item(10,13).Value = 0
$ws5.Cells.Item(10,14).Value = 0
$ws5.Cells.Item(11,1).Value = 9
$ws5.Cells.Item(11,2).Value = 0
$ws5.Cells.Item(11,3).Value = 0
$ws5.Cells.Item(11,4).Value = 0
$ws5.Cells.Item(11,5).Value = 0
$ws5.Cells.Item(11,6).Value = 0
$ws5.Cells.Item(11,7).Value = 402.5167240182
$ws5.Cells.Item(11,8).Value = 391.1185967199499
$ws5.Cells.Item(11,9).Value = 2755.090198377
$ws5.Cells.Item(11,10).Value = 124.72807872085
$ws5.Cells.Item(11,11).Value = 0
$ws5.Cells.Item(11,12).Value = 0
$ws5.Cells.Item(11,13).Value = 0
$ws5.Cells.Item(11,14).Value = 0
$ws5.Cells.Item(12,1).Value = 10
$ws5.Cells.Item(12,2).Value = 0
$ws5.Cells.Item(12,3).Value = 0
$ws5.Cells.Item(12,4).Value = 0
$ws5.Cells.Item(12,5).Value = 0
$ws5.Cells.Item(12,6).Value = 0
$ws5.Cells.Item(12,7).Value = 1487.315532974925
$ws5.Cells.Item(12,8).Value = 1299.814278356517
$ws5.Cells.Item(12,9).Value = 211.003241702325
$ws5.Cells.Item(12,10).Value = 273.78957216009
$ws5.Cells.Item(12,11).Value = 0
$ws5.Cells.Item(12,12).Value = 0
$ws5.Cells.Item(12,13).Value = 0
$ws5.Cells.Item(12,14).Value = 0
$ws5.Cells.Item(13,1).Value = 11
$ws5.Cells.Item(13,2).Value = 0
$ws5.Cells.Item(13,3).Value = 0
$ws5.Cells.Item(13,4).Value = 0
$ws5.Cells.Item(13,5).Value = 0
$ws5.Cells.Item(13,6).Value = 0
$ws5.Cells.Item(13,7).Value = 1369.662622241891
$ws5.Cells.Item(13,8).Value = 1103.802765508102
$ws5.Cells.Item(13,9).Value = 2440.815931659878
$ws5.Cells.Item(13,10).Value = 0
$ws5.Cells.Item(13,11).Value = 0
$ws5.Cells.Item(13,12).Value = 0
$ws5.Cells.Item(13,13).Value = 0
$ws5.Cells.Item(13,14).Value = 0
$ws5.Cells.Item(14,1).Value = 12
$ws5.Cells.Item(14,2).Value = 0
$ws5.Cells.Item(14,3).Value = 0
$ws5.Cells.Item(14,4).Value = 0
$ws5.Cells.Item(14,5).Value = 0
$ws5.Cells.Item(14,6).Value = 0
$ws5.Cells.Item(14,7).Value = 132.7547882557438
$ws5.Cells.Item(14,8).Value = 0
$ws5.Cells.Item(14,9).Value = 2095.549248885259
$ws5.Cells.Item(14,10).Value = 275.950964351827
$ws5.Cells.Item(14,11).Value = 0
$ws5.Cells.Item(14,12).Value = 0
$ws5.Cells.Item(14,13).Value = 0
$ws5.Cells.Item(14,14).Value = 0
$ws5.Cells.Item(15,1).Value = 13
$ws5.Cells.Item(15,2).Value = 0
$ws5.Cells.Item(15,3).Value = 0
$ws5.Cells.Item(15,4).Value = 0
$ws5.Cells.Item(15,5).Value = 0
$ws5.Cells.Item(15,6).Value = 0
$ws5.Cells.Item(15,7).Value = 631.8318778127319
$ws5.Cells.Item(15,8).Value = 0
$ws5.Cells.Item(15,9).Value = 540.6336848907233
$ws5.Cells.Item(15,10).Value = 531.186091966456
$ws5.Cells.Item(15,11).Value = 0
$ws5.Cells.Item(15,12).Value = 0
$ws5.Cells.Item(15,13).Value = 0
$ws5.Cells.Item(15,14).Value = 0
$ws5.Cells.Item(16,1).Value = 14
$ws5.Cells.Item(16,2).Value = 0
$ws5.Cells.Item(16,3).Value = 0
$ws5.Cells.Item(16,4).Value = 0
$ws5.Cells.Item(16,5).Value = 0
$ws5.Cells.Item(16,6).Value = 0
$ws5.Cells.Item(16,7).Value = 0
$ws5.Cells.Item(16,8).Value = 9.529306414864212
$ws5.Cells.Item(16,9).Value = 2833.380440686292
$ws5.Cells.Item(16,10).Value = 0
$ws5.Cells.Item(16,11).Value = 0
$ws5.Cells.Item(16,12).Value = 0
$ws5.Cells.Item(16,13).Value = 0
$ws5.Cells.Item(16,14).Value = 0
$ws5.Cells.Item(17,1).Value = 15
$ws5.Cells.Item(17,2).Value = 0
$ws5.Cells.Item(17,3).Value = 0
$ws5.Cells.Item(17,4).Value = 0
$ws5.Cells.Item(17,5).Value = 0
$ws5.Cells.Item(17,6).Value = 0
$ws5.Cells.Item(17,7).Value = 0
$ws5.Cells.Item(17,8).Value = 721.4150622771695
$ws5.Cells.Item(17,9).Value = 0
$ws5.Cells.Item(17,10).Value = 0
$ws5.Cells.Item(17,11).Value = 0
$ws5.Cells.Item(17,12).Value = 0
$ws5.Cells.Item(17,13).Value = 0
$ws5.Cells.Item(17,14).Value = 0

# --- Capacity ---
$ws6 = $wb.Worksheets.Item("Capacity")
$ws6.Cells.Item(2,2).Value = 793720
$ws6.Cells.Item(3,2).Value = 115917.5
$ws6.Cells.Item(4,2).Value = 112217.5
$ws6.Cells.Item(5,2).Value = 132746.25
$ws6.Cells.Item(6,2).Value = 75023.125
$ws6.Cells.Item(7,2).Value = 28979.375
$ws6.Cells.Item(8,2).Value = 112217.5
$ws6.Cells.Item(9,2).Value = 331865.625
$ws6.Cells.Item(10,2).Value = 43250
$ws6.Cells.Item(11,2).Value = 115917.5
$ws6.Cells.Item(12,2).Value = 84163.125
$ws6.Cells.Item(13,2).Value = 199119.375
$ws6.Cells.Item(14,2).Value = 75023.125

# --- ProcessingTime ---
$ws7 = $wb.Worksheets.Item("ProcessingTime")
$ws7.Cells.Item(2,2).Value = 4
$ws7.Cells.Item(3,3).Value = 4
$ws7.Cells.Item(4,4).Value = 4
$ws7.Cells.Item(5,5).Value = 2
$ws7.Cells.Item(6,6).Value = 1
$ws7.Cells.Item(7,7).Value = 1
$ws7.Cells.Item(8,8).Value = 4
$ws7.Cells.Item(9,9).Value = 5
$ws7.Cells.Item(10,10).Value = 5
$ws7.Cells.Item(11,11).Value = 4
$ws7.Cells.Item(12,12).Value = 3
$ws7.Cells.Item(13,13).Value = 3

